# "Lunchtimes and lecturer time limits implemented."
# Bump the per-subject lecture-hours row (row 2) from 2 hours to 7 hours
# for every subject column (A:AN), and move the viewport/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lecture hours for all 40 subjects: 2 -> 7
$ws.Range("A2:AN2").Value = 7

# Scroll the window so column W is left-most visible, and select AL6
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 23
$ws.Range("AL6").Select()
